$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.090.75'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.94%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.891.29'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +1.82%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9990'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '307.67'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +0.92%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9983'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5187'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +2.16%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3726'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +2.04%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07205'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +0.42%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.9053'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +1.77%  '
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +1.95%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.07655'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +1.64%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.878.59'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +1.18%  '
$ws.Range('B14').Value = 'Litecoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '95.17'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +3.74%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.279'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +0.76%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.9982'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -0.25%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.000008501'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -0.26%  '
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +1.85%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.9985'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -0.15%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '27.140.40'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.053'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +0.79%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '2.104.06'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +0.46%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '10.56'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +2.42%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.463'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '145.76'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -0.19%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.790'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -1.42%  '
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +1.29%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.150'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +5.04%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '114.65'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +1.58%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.939'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +5.65%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.800'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +3.43%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.09206'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -0.25%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.05048'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -0.92%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7624'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +3.95%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.195'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +4.08%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.017'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.280'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +2.33%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.557'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +4.10%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.5627'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +6.01%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.01996'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -0.61%  '
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +0.31%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '6.596'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +1.51%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '118.46'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +0.24%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.863'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +5.59%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.4803'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +3.69%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.15'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +2.66%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.9979'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -0.18%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.577'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +1.17%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '37.12'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +0.44%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '63.52'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +0.93%  '
